$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AA2 with new multi-line time log text
$ws.Range("AA2").Value = "14:54-14:55(No thing)`n14:55-14:59(No thing)`n15:00-15:01(No thing)`n15:02-15:02(No thing)`n"

# Set AG2 (Tong/Total) to computed total hours 0.1 (stored as text)
$ws.Range("AG2").Value = "'0.1"

# Clear AA6 (time log entries removed after being processed)
$ws.Range("AA6").Value = ""
